# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp stamps that get refreshed each time
# the handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 84f75f2a-...md (row 3, col G)
$wsOverview.Range("G3").Value = "2016-08-13 17:02:15"

# zh-cn sheet: Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
# for 84f75f2a-...9978baf4...zh-cn.xlf (row 3)
$wsZhCn.Range("H3").Value = "2016-08-13 17:02:07"
$wsZhCn.Range("K3").Value = "2016-08-13 17:02:38"

# de-de sheet: Correspond Handoff Datetime (H3) for 84f75f2a-...md (row 3) —
# this cell shares its value with the Overview sheet's "Latest HO Xliff
# Generate Date" stamp, so it must be refreshed too.
$wsDeDe.Range("H3").Value = "2016-08-13 17:02:15"

# de-de sheet: Correspond Handback DateTime (K3) for
# 84f75f2a-...9978baf4...de-de.xlf (row 3)
$wsDeDe.Range("K3").Value = "2016-08-13 17:02:48"
